# Apply the "RF diameter (pixels)" column rework to the regions.csv sheet.
#
# Summary of the change (see commit "more rough work for neuroscience club
# presentation"):
#   - a new "RF diameter (pixels)" shared string is introduced and the old
#     "RF diameter" string is retired
#   - the F..L columns of the data table are reshuffled:
#       new F = old G ("RF diameter", now literal pixel counts, header
#                       renamed to "RF diameter (pixels)")
#       new G = FLOOR(new F / 2, 1)           (was old H, "RF Stride")
#       new H = ROUND(new L / I / J, 0)        (was old F, "# features")
#       new K = old L (literal "order" value)
#       new L = old K ("pixels * features" formula, now keyed off H)
#   - the selected cell moves to F12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("regions.csv")

# ---------------------------------------------------------------------
# Header row (row 1) - relabel columns F, G, H, K, L
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "RF diameter (pixels)"
$ws.Range("G1").Value = "RF Stride"
$ws.Range("H1").Value = "# features"
$ws.Range("K1").Value = "order (to distinguish ff and fb connections)"
$ws.Range("L1").Value = "pixels * features"

# ---------------------------------------------------------------------
# Row 2 - the "base" row: literal RF diameter / #features, derived stride
# ---------------------------------------------------------------------
$ws.Range("F2").Value = 5
$ws.Range("G2").Formula = '=FLOOR(F2/2,1)'
$ws.Range("H2").Value = 364
$ws.Range("I2").Formula = '=360*D2/G2'
$ws.Range("J2").Formula = '=640*D2/G2'
$ws.Range("K2").Value = 1
$ws.Range("L2").Formula = '=H2*I2*J2'

# ---------------------------------------------------------------------
# Rows 3-17, 19, 21 - regular formula rows
#   F = C*F$2 ; G = FLOOR(F/2,1) ; H = ROUND(L/I/J,0)
#   I = 360*D/G ; J = 640*D/G ; K = literal ; L = L$2*B/B$2
# ---------------------------------------------------------------------
$orderValues = @{
    3  = 2
    4  = 3
    5  = 4
    6  = 3
    7  = 5
    8  = 6
    9  = 6
    10 = 6
    11 = 6
    12 = 6
    13 = 7
    14 = 6
    15 = 7
    16 = 7
    17 = 6
    19 = 7
    21 = 6
}

foreach ($r in @(3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,19,21)) {
    $ws.Range("F$r").Formula = "=C$r*F`$2"
    $ws.Range("G$r").Style = "Normal"
    $ws.Range("G$r").Formula = "=FLOOR(F$r/2,1)"
    $ws.Range("H$r").Formula = "=ROUND(L$r/I$r/J$r,0)"
    $ws.Range("I$r").Formula = "=360*D$r/G$r"
    $ws.Range("J$r").Formula = "=640*D$r/G$r"
    $ws.Range("K$r").Value = $orderValues[$r]
    $ws.Range("L$r").Formula = "=L`$2*B$r/B`$2"
}

# ---------------------------------------------------------------------
# Rows 18 & 20 - the two "output" rows (literal RF diameter, derived
# stride, literal #features/order)
# ---------------------------------------------------------------------
foreach ($r in @(18,20)) {
    $ws.Range("H$r").Value = 1
    $ws.Range("K$r").Value = 8
}

# row 18's literal RF-diameter cell keeps the highlighted manual-entry
# style that the old "RF diameter" column (G13:G21) used
$ws.Range("D7").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("F18").Value = 100
$ws.Range("G18").Formula = '=FLOOR(F18/2,1)'

$ws.Range("F20").Value = 100
$ws.Range("G20").Formula = '=FLOOR(F20/2,1)'

# ---------------------------------------------------------------------
# Move the active selection to F12 (matches the author's last click)
# ---------------------------------------------------------------------
$ws.Range("F12").Select()

Write-Output "RF diameter (pixels) rework applied"
